## ---------------------------------------------------------------------
## OSI-ZADACI-ZA-CLANOVE-TIMA.docx - "Doradjena korisnicka dokumentacija,
## klasa Datum i testirano"
##
## Four content edits, applied by rebuilding the affected paragraphs'
## OOXML in-place (keeps every paragraph/run formatting property exactly
## as Word would emit it, including <w:proofErr/> spell-check markers).
## ---------------------------------------------------------------------

$d = $word.ActiveDocument

function Find-ParagraphIndexContaining([string]$needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -match [regex]::Escape($needle)) {
            return $i
        }
    }
    return -1
}

$W_OPEN = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/rewrite.xml" pkg:contentType="text/xml"><pkg:xmlData>'
$W_CLOSE = '</pkg:xmlData></pkg:part></pkg:package>'

## ======================================================================
## Edit 1: "- .тхт фајл ... од (tip klase Datum), до (tip klase Datum),
## вријеме (tip klase datum), ..." paragraph -> first two "(tip klase
## Datum)" become "(sati, minuti, sekunde)" / "(sati, minute, sekunde)";
## the third, lower-case, "(tip klase datum)" is left untouched.
## ======================================================================

$idxKlase = Find-ParagraphIndexContaining("klase")

$xmlKlase = @"
$W_OPEN
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:ind w:left="720"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="sr-Cyrl-RS"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="sr-Cyrl-RS"/>
    </w:rPr>
    <w:t xml:space="preserve">- .тхт фајл ће имати у себи колоне </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="sr-Cyrl-RS"/>
    </w:rPr>
    <w:t xml:space="preserve">ознака лета, </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="sr-Cyrl-RS"/>
    </w:rPr>
    <w:t>од</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> (</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve">sati, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>minuti</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve">, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>sekunde</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>)</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="sr-Cyrl-RS"/>
    </w:rPr>
    <w:t>, до</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> (</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve">sati, minute, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>sekunde</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>)</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="sr-Cyrl-RS"/>
    </w:rPr>
    <w:t>, вријеме</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> (tip </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>klase</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> datum)</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="sr-Cyrl-RS"/>
    </w:rPr>
    <w:t>, опис, број мјеста, број тренутно слободних мјеста</w:t>
  </w:r>
</w:p>
$W_CLOSE
"@

$d.Paragraphs.Item($idxKlase).Range.InsertXML($xmlKlase)

## ======================================================================
## Edit 2: new bullet "Smisliti kako izgledaju te rezervacije" inserted
## between "Омогућити отварање..." and "Ако не постоји лет...";
## the _GoBack bookmark moves to the end of the new bullet, and
## "Ако не постоји..." gains a <w:lastRenderedPageBreak/> before its text.
## ======================================================================

$idxOmoguciti = Find-ParagraphIndexContaining("Омогућити отварање")
$pOmoguciti = $d.Paragraphs.Item($idxOmoguciti)
$pOmoguciti.Range.InsertParagraphAfter() | Out-Null

$idxNew = $idxOmoguciti + 1
$pNew = $d.Paragraphs.Item($idxNew)

$xmlNewBullet = @"
$W_OPEN
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="sr-Cyrl-RS"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>Smisliti</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>kako</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>izgledaju</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>te</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> rezervacije</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
$W_CLOSE
"@

$pNew.Range.InsertXML($xmlNewBullet)

$idxAko = Find-ParagraphIndexContaining("Ако не постоји лет")
$pAko = $d.Paragraphs.Item($idxAko)

$xmlAko = @"
$W_OPEN
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="sr-Cyrl-RS"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
      <w:lang w:val="sr-Cyrl-RS"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t xml:space="preserve">Ако не постоји лет, то јест фајл за резервацију лета враћа се одговарајућа порука раднику </w:t>
  </w:r>
</w:p>
$W_CLOSE
"@

$pAko.Range.InsertXML($xmlAko)

## ======================================================================
## Edit 3: final paragraph "... iz opravdanih je razloga." -> the trailing
## text is re-split with a <w:proofErr/> pair around "razloga", and the
## _GoBack bookmark (now living on the new bullet from edit 2) is removed
## from here.
## ======================================================================

$idxRazloga = Find-ParagraphIndexContaining("opravdanih")
$pRazloga = $d.Paragraphs.Item($idxRazloga)

$xmlRazloga = @"
$W_OPEN
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>Izvinjavam</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> se </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:proofErr w:type="gramStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>na</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:proofErr w:type="gramEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>miksu</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>cirilice</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>i</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>latinice</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve">, </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>iz</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>opravdanih</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t xml:space="preserve"> je </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>razloga</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:sz w:val="28"/>
      <w:szCs w:val="28"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
</w:p>
$W_CLOSE
"@

$pRazloga.Range.InsertXML($xmlRazloga)

Write-Host "Edits applied."
